$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column after column A (old B,C,D -> C,D,E)
$ws.Columns.Item(2).Insert()

# 2) Insert two new rows above row 1 (everything shifts down by 2)
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# --- at this point the old sheet (A1:D17) now lives at A3:E19 ---

# 3) New title row (row 1): "Erhebungstools (Engere Auswahl)" in column B, bold 14pt
$title = $ws.Range("B1")
$title.Value = "Erhebungstools (Engere Auswahl)"
$title.Font.Bold = $true
$title.Font.Size = 14
$ws.Rows.Item(1).RowHeight = 18.75

# 4) Row 2 stays empty, just a thin bottom divider
$ws.Rows.Item(2).RowHeight = 15.75

# 5) Row 3 (former row 1): label "Anbieter:" in A3:B3, bold 14pt like other section headers
$a3 = $ws.Range("A3:B3")
$a3.Font.Bold = $true
$a3.Font.Size = 14
$ws.Range("A3").Value = "Anbieter:"
$ws.Rows.Item(3).RowHeight = 19.5

# 6) The two merged section-header rows (former rows 2 and 7, now rows 4 and 9)
#    get the same bold 14pt treatment across their full width
$hdr1 = $ws.Range("A4:E4")
$hdr1.Font.Bold = $true
$hdr1.Font.Size = 14
$ws.Rows.Item(4).RowHeight = 19.5

$hdr2 = $ws.Range("A9:E9")
$hdr2.Font.Bold = $true
$hdr2.Font.Size = 14
$ws.Rows.Item(9).RowHeight = 19.5

# 7) Column widths: new column B plus restored widths for the shifted columns
$ws.Columns.Item(2).ColumnWidth = 27.75
